$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.068.60"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "2.522.93"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.43"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.09"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").Value = "2.515.85"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("E11").Value = "  -2.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "2.956.63"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.97"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "59.001.71"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").Value = "2.503.74"
$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.06"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.26"
$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.49"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.15"
$ws.Range("E31").Value = "  +3.46%  "

$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +6.22%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.36"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  -2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.08"
$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.17"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.39"
$ws.Range("E46").Value = "  +4.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.85"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0922"
$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0219"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.31"
$ws.Range("E51").Value = "  -1.82%  "
